# Edit slide 9 ("STAMP - Destination Node Address TLV"), shape "Rectangle 8":
#  - append " (e.g. sweeping ECMP paths)" before the trailing period of the
#    last bullet's text
#  - the shape has spAutoFit, so its height grows to fit the extra text

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$sh = $s.Shapes.Item(5)

$para = $sh.TextFrame.TextRange.Paragraphs(5)
$run = $para.Runs(1)
$run.Text = "Useful when query is sent with 127/8 destination address (e.g. sweeping ECMP paths)."

# Match the exact auto-fit height PowerPoint computed for the new text.
$sh.Height = 2772234 / 12700
